# Sprint_Backlog_3.xlsx — "Add files via upload" re-upload.
#
# Net effect versus the previous revision: the task row for "HU-3 T-4"
# ("Implementar la asociación entre un objetivo y los departamentos
# seleccionados con sus validaciones correspondientes") was removed from
# the Tasks sheet. Deleting the whole row shifts every following row up
# by one (old row 10 -> new row 9, ... old row 17 -> new row 16) and lets
# Excel drop the now-unused shared strings ("Implementar la asociación…"
# and "HU-3 T-4") automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 9 holds the "HU-3 T-4" task (A9/B9 etc.); deleting it shifts the
# rows below (old 10..17) up into 9..16.
$ws.Rows(9).Delete()

# Match the author's final selection/active cell after the edit.
$ws.Range("A8").Select()
